# "deletion of billing in client side"
# The client re-billed the invoice under a new company/address, refreshed
# the two timestamped line entries, updated the first line's amount, and
# removed (deleted) the second billing line entirely (row 18) along with
# its trailing SPO note (row 19).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: client name / address -----------------------------------
$ws.Range("B9").Value = "COMPLEX TECH"
$ws.Range("A10").Value = "COMPLEX, STA. ROSA, LAGUNA"

# --- Refresh the "printed at" timestamp (H9) --------------------------
# Evaluate NOW() to get a plain numeric date serial (so the cell keeps
# storing a literal value, no formula - matching how the sheet already
# stores these timestamps).
$now = $excel.Evaluate("NOW()")

$ws.Range("H9").Value = $now

# --- Line 1 (row 16): timestamp refreshed, amount updated -------------
$ws.Range("B16").Value = $now
$ws.Range("H16").Value = 565

# --- Line 2 (row 18) deleted from the client side, along with its -----
# --- trailing SPO reference row (row 19) -------------------------------
$ws.Range("A18:H18").ClearContents()
$ws.Range("D19:E19").ClearContents()
